$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
$ws.Range("F2").Value = 31
$ws.Range("G2").Value = "adam"
$ws.Range("I2").Value = 64
$ws.Range("J2").Value = 29.06751262703127
$ws.Range("K2").Value = 1382.199652640461
$ws.Range("L2").Value = 37.17794578295661
$ws.Range("M2").Value = 0.184509935589258

# --- Row 3 (new) ---
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "LSTM"
$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 60
$ws.Range("F3").Value = 31
$ws.Range("G3").Value = "<keras.src.optimizers.legacy.adam.Adam object at 0x7f62c62e4940>"
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 32
$ws.Range("J3").Value = 20.46878326057768
$ws.Range("K3").Value = 816.5461207380611
$ws.Range("L3").Value = 28.57527114023699
$ws.Range("M3").Value = 0.1288886322929146

# --- Row 4 (new) ---
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "LSTM"
$ws.Range("D4").Value = 40
$ws.Range("E4").Value = 60
$ws.Range("F4").Value = 31
$ws.Range("G4").Value = "adam"
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 64
$ws.Range("J4").Value = 34.00734830684019
$ws.Range("K4").Value = 2007.035158362287
$ws.Range("L4").Value = 44.79994596383223
$ws.Range("M4").Value = 0.2136256006984066

# --- Row 5 (new) ---
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "LSTM"
$ws.Range("D5").Value = 40
$ws.Range("E5").Value = 60
$ws.Range("F5").Value = 31
$ws.Range("G5").Value = "<keras.src.optimizers.legacy.adam.Adam object at 0x7f624bebb5b0>"
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 32
$ws.Range("J5").Value = 48.06092926304213
$ws.Range("K5").Value = 3699.469716546063
$ws.Range("L5").Value = 60.82326624365106
$ws.Range("M5").Value = 0.2958491047321462
